# Regenerate the "K" column (G) values for save_data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6
$ws.Range("G3").Value = 4
$ws.Range("G4").Value = 3
$ws.Range("G5").Value = 4
$ws.Range("G6").Value = 5
$ws.Range("G7").Value = 5
$ws.Range("G8").Value = 4
